$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: column A takes column B's old width, and vice versa.
$ws.Columns.Item(1).ColumnWidth = 16.42578125
$ws.Columns.Item(2).ColumnWidth = 15.7109375

# Update cell values
$ws.Range("A1").Value = 0.019719031323889441
$ws.Range("B1").Value = -0.019719031381164393

$ws.Range("A2").Value = -0.00079887094455772579
$ws.Range("B2").Value = 0.00079887083129568336

$ws.Range("A3").Value = -0.035171490645430063
$ws.Range("B3").Value = 0.035171490585610032

$ws.Range("A4").Value = -0.04070507555833966
$ws.Range("B4").Value = 0.040705075470485512

$ws.Range("A5").Value = 0.034173019087727019
$ws.Range("B5").Value = -0.034173019229503283
